# TaskBreakDown sheet update (AashishSharma.xlsx)
# - Hours Burnt (F) / Remaining Hours (G, formula =E-F) updates for rows 27-30
# - Window/selection moved down to F29

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: Hours Burnt 0 -> 3 (Remaining Hours recalculates via existing formula)
$ws.Range("F27").Value = 3

# Row 28: Hours Burnt 1 -> 4
$ws.Range("F28").Value = 4

# Row 29: Hours Burnt 1 -> 5
$ws.Range("F29").Value = 5

# Row 30: Hours Burnt 1 -> 4
$ws.Range("F30").Value = 4

# Move the active selection / viewport down to F29 (matches the saved view state)
$ws.Range("F29").Select()

# Scroll the window so row 19 is at the top-left (topLeftCell A14 -> A19)
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
